$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in cell A11: "Mentaland behavioraldisorders" -> "Mental and behavioral disorders"
$ws.Range("A11").Value = "F                                                    Mental and behavioral disorders                                  88                            135                         0.59 (0.54-0.64)                 0.65 (0.59-0.70)  "

# Update the active selection to A12
$ws.Range("A12").Select()
